# Insert two new data rows at the top of the data block (row 247), pushing
# the existing rows (247-314) down to (249-316). This represents the weekly
# refresh described in the commit message: a new week's worth of "Apio"
# price data for Feria Lagunitas de Puerto Montt was published, and it is
# prepended to the sheet's data (directly below the header + the single
# already-present block), while the previously-most-recent rows shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows starting at row 247. This shifts the existing rows
# 247-314 down to 249-316, matching the rest of the sheet's existing
# formatting since Excel's Insert copies formatting from the row above.
$ws.Rows.Item(247).Resize(2).Insert()

# New row 247: "Primera" quality entry for the new date (2022-09-05).
$ws.Cells.Item(247, 1).Value = 4
$ws.Cells.Item(247, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(247, 3).Value = "Los Lagos"
$ws.Cells.Item(247, 4).Value = 44809
$ws.Cells.Item(247, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(247, 5).Value = 10
$ws.Cells.Item(247, 6).Value = 100112017
$ws.Cells.Item(247, 7).Value = "Apio"
$ws.Cells.Item(247, 8).Value = "Americana (o)"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 15
$ws.Cells.Item(247, 11).Value = 15000
$ws.Cells.Item(247, 12).Value = 15000
$ws.Cells.Item(247, 13).Value = 15000
$ws.Cells.Item(247, 14).Value = "$/docena de matas"
$ws.Cells.Item(247, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(247, 16).Value = 2500
$ws.Cells.Item(247, 17).Value = 6
$ws.Cells.Item(247, 18).Value = "Hortaliza"

# New row 248: "Segunda" quality entry for the same new date.
$ws.Cells.Item(248, 1).Value = 4
$ws.Cells.Item(248, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(248, 3).Value = "Los Lagos"
$ws.Cells.Item(248, 4).Value = 44809
$ws.Cells.Item(248, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(248, 5).Value = 10
$ws.Cells.Item(248, 6).Value = 100112017
$ws.Cells.Item(248, 7).Value = "Apio"
$ws.Cells.Item(248, 8).Value = "Americana (o)"
$ws.Cells.Item(248, 9).Value = "Segunda"
$ws.Cells.Item(248, 10).Value = 15
$ws.Cells.Item(248, 11).Value = 12000
$ws.Cells.Item(248, 12).Value = 12000
$ws.Cells.Item(248, 13).Value = 12000
$ws.Cells.Item(248, 14).Value = "$/docena de matas"
$ws.Cells.Item(248, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(248, 16).Value = 2000
$ws.Cells.Item(248, 17).Value = 6
$ws.Cells.Item(248, 18).Value = "Hortaliza"
